$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update order number and its date in columns A and B for row 4
$ws.Range("A4").Value = "20240716-1"
$ws.Range("B4").Value = "16.07.2024"

# Set the readiness date in column X for row 4
$ws.Range("X4").Value = "16.07.2024"
